$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 1287, shifting the existing 1287:1373 block down to 1291:1377
$ws.Rows("1287:1290").Insert()

# ---- New row 1287 ----
$ws.Cells.Item(1287, 1).Value = 11
$ws.Cells.Item(1287, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(1287, 3).Value = "Bíobío"
$ws.Cells.Item(1287, 4).Value = 45223
$ws.Cells.Item(1287, 5).Value = 8
$ws.Cells.Item(1287, 6).Value = 100112033
$ws.Cells.Item(1287, 7).Value = "Lechuga"
$ws.Cells.Item(1287, 8).Value = "Conconina(o)"
$ws.Cells.Item(1287, 9).Value = "Primera"
$ws.Cells.Item(1287, 10).Value = 120
$ws.Cells.Item(1287, 11).Value = 10000
$ws.Cells.Item(1287, 12).Value = 10000
$ws.Cells.Item(1287, 13).Value = 10000
$ws.Cells.Item(1287, 14).Value = "$/caja 10 unidades"
$ws.Cells.Item(1287, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1287, 16).Value = 1000
$ws.Cells.Item(1287, 17).Value = 10
$ws.Cells.Item(1287, 18).Value = "Hortaliza"

# ---- New row 1288 ----
$ws.Cells.Item(1288, 1).Value = 11
$ws.Cells.Item(1288, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(1288, 3).Value = "Bíobío"
$ws.Cells.Item(1288, 4).Value = 45223
$ws.Cells.Item(1288, 5).Value = 8
$ws.Cells.Item(1288, 6).Value = 100112033
$ws.Cells.Item(1288, 7).Value = "Lechuga"
$ws.Cells.Item(1288, 8).Value = "Conconina(o)"
$ws.Cells.Item(1288, 9).Value = "Segunda"
$ws.Cells.Item(1288, 10).Value = 150
$ws.Cells.Item(1288, 11).Value = 8000
$ws.Cells.Item(1288, 12).Value = 8000
$ws.Cells.Item(1288, 13).Value = 8000
$ws.Cells.Item(1288, 14).Value = "$/caja 12 unidades"
$ws.Cells.Item(1288, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1288, 16).Value = 667
$ws.Cells.Item(1288, 17).Value = 12
$ws.Cells.Item(1288, 18).Value = "Hortaliza"

# ---- New row 1289 ----
$ws.Cells.Item(1289, 1).Value = 11
$ws.Cells.Item(1289, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(1289, 3).Value = "Bíobío"
$ws.Cells.Item(1289, 4).Value = 45223
$ws.Cells.Item(1289, 5).Value = 8
$ws.Cells.Item(1289, 6).Value = 100112033
$ws.Cells.Item(1289, 7).Value = "Lechuga"
$ws.Cells.Item(1289, 8).Value = "Escarola"
$ws.Cells.Item(1289, 9).Value = "Primera"
$ws.Cells.Item(1289, 10).Value = 150
$ws.Cells.Item(1289, 11).Value = 13000
$ws.Cells.Item(1289, 12).Value = 13000
$ws.Cells.Item(1289, 13).Value = 13000
$ws.Cells.Item(1289, 14).Value = "$/caja 15 unidades"
$ws.Cells.Item(1289, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(1289, 16).Value = 867
$ws.Cells.Item(1289, 17).Value = 15
$ws.Cells.Item(1289, 18).Value = "Hortaliza"

# ---- New row 1290 ----
$ws.Cells.Item(1290, 1).Value = 11
$ws.Cells.Item(1290, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(1290, 3).Value = "Bíobío"
$ws.Cells.Item(1290, 4).Value = 45223
$ws.Cells.Item(1290, 5).Value = 8
$ws.Cells.Item(1290, 6).Value = 100112033
$ws.Cells.Item(1290, 7).Value = "Lechuga"
$ws.Cells.Item(1290, 8).Value = "Escarola"
$ws.Cells.Item(1290, 9).Value = "Segunda"
$ws.Cells.Item(1290, 10).Value = 120
$ws.Cells.Item(1290, 11).Value = 10000
$ws.Cells.Item(1290, 12).Value = 10000
$ws.Cells.Item(1290, 13).Value = 10000
$ws.Cells.Item(1290, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(1290, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(1290, 16).Value = 556
$ws.Cells.Item(1290, 17).Value = 18
$ws.Cells.Item(1290, 18).Value = "Hortaliza"
